$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "eng_lang"
$ws1.Range("D7").Value = "Чоп этилган материаллар номи (инглиз тилида)"
$ws1.Range("D7").Select()

$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "rus_uzb_lang"

